# Mise à jour de l'application
# Append a new GPS training session ("Entrainement", J-3, 2025-09-03)
# as 13 new rows at the bottom of the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 417
$newRowCount = 13
$templateRow = $firstNewRow - 1   # row 416: last existing data row, used as a style/format template

# Insert blank rows by copying the last existing data row repeatedly, so the
# date number format (column B) and the centered style (column D) carry over
# to the new rows exactly like they do for the existing rows above.
for ($i = 0; $i -lt $newRowCount; $i++) {
    $ws.Rows.Item($templateRow).Copy()
    $ws.Rows.Item($firstNewRow).Insert()
}

# Row 417: Mattheo Haon (right back)
$ws.Range("A417").Value = 'Entrainement'
$ws.Range("B417").Value = 45903
$ws.Range("C417").Value = 'Global'
$ws.Range("D417").Value = 'J-3'
$ws.Range("E417").Value = 'Mattheo Haon'
$ws.Range("F417").Value = 'right back'
$ws.Range("G417").Value = '01:26:49'
$ws.Range("H417").Value = 7.06
$ws.Range("I417").Value = 0.99
$ws.Range("J417").Value = 6.05
$ws.Range("K417").Value = 0.59
$ws.Range("L417").Value = 0.3
$ws.Range("M417").Value = 0.12
$ws.Range("N417").Value = 0
$ws.Range("O417").Value = 8
$ws.Range("P417").Value = 4.85
$ws.Range("Q417").Value = 30.08
$ws.Range("R417").Value = 4.22
$ws.Range("S417").Value = 31
$ws.Range("T417").Value = 3
$ws.Range("U417").Value = 12
$ws.Range("V417").Value = 10

# Row 418: Romain Thunet (center back)
$ws.Range("A418").Value = 'Entrainement'
$ws.Range("B418").Value = 45903
$ws.Range("C418").Value = 'Global'
$ws.Range("D418").Value = 'J-3'
$ws.Range("E418").Value = 'Romain Thunet'
$ws.Range("F418").Value = 'center back'
$ws.Range("G418").Value = '01:26:03'
$ws.Range("H418").Value = 6.23
$ws.Range("I418").Value = 0.5
$ws.Range("J418").Value = 5.72
$ws.Range("K418").Value = 0.39
$ws.Range("L418").Value = 0.11
$ws.Range("M418").Value = 0.01
$ws.Range("N418").Value = 0
$ws.Range("O418").Value = 2
$ws.Range("P418").Value = 4.22
$ws.Range("Q418").Value = 25.21
$ws.Range("R418").Value = 4.05
$ws.Range("S418").Value = 19
$ws.Range("T418").Value = 2
$ws.Range("U418").Value = 14
$ws.Range("V418").Value = 2

# Row 419: Ilyes Boughanmi (center forward)
$ws.Range("A419").Value = 'Entrainement'
$ws.Range("B419").Value = 45903
$ws.Range("C419").Value = 'Global'
$ws.Range("D419").Value = 'J-3'
$ws.Range("E419").Value = 'Ilyes Boughanmi'
$ws.Range("F419").Value = 'center forward'
$ws.Range("G419").Value = '01:29:06'
$ws.Range("H419").Value = 6.31
$ws.Range("I419").Value = 0.8
$ws.Range("J419").Value = 5.5
$ws.Range("K419").Value = 0.39
$ws.Range("L419").Value = 0.34
$ws.Range("M419").Value = 0.09
$ws.Range("N419").Value = 0
$ws.Range("O419").Value = 10
$ws.Range("P419").Value = 4.18
$ws.Range("Q419").Value = 29.46
$ws.Range("R419").Value = 5.05
$ws.Range("S419").Value = 25
$ws.Range("T419").Value = 5
$ws.Range("U419").Value = 14
$ws.Range("V419").Value = 4

# Row 420: Jeremie Laurent (left forward)
$ws.Range("A420").Value = 'Entrainement'
$ws.Range("B420").Value = 45903
$ws.Range("C420").Value = 'Global'
$ws.Range("D420").Value = 'J-3'
$ws.Range("E420").Value = 'Jeremie Laurent'
$ws.Range("F420").Value = 'left forward'
$ws.Range("G420").Value = '01:26:26'
$ws.Range("H420").Value = 6.61
$ws.Range("I420").Value = 0.98
$ws.Range("J420").Value = 5.62
$ws.Range("K420").Value = 0.52
$ws.Range("L420").Value = 0.38
$ws.Range("M420").Value = 0.1
$ws.Range("N420").Value = 0
$ws.Range("O420").Value = 8
$ws.Range("P420").Value = 4.55
$ws.Range("Q420").Value = 28.66
$ws.Range("R420").Value = 5.06
$ws.Range("S420").Value = 30
$ws.Range("T420").Value = 12
$ws.Range("U420").Value = 21
$ws.Range("V420").Value = 6

# Row 421: Omar Benyounes (center midfield)
$ws.Range("A421").Value = 'Entrainement'
$ws.Range("B421").Value = 45903
$ws.Range("C421").Value = 'Global'
$ws.Range("D421").Value = 'J-3'
$ws.Range("E421").Value = 'Omar Benyounes'
$ws.Range("F421").Value = 'center midfield'
$ws.Range("G421").Value = '01:27:42'
$ws.Range("H421").Value = 6.5
$ws.Range("I421").Value = 0.87
$ws.Range("J421").Value = 5.62
$ws.Range("K421").Value = 0.47
$ws.Range("L421").Value = 0.31
$ws.Range("M421").Value = 0.09
$ws.Range("N421").Value = 0
$ws.Range("O421").Value = 9
$ws.Range("P421").Value = 4.37
$ws.Range("Q421").Value = 28.95
$ws.Range("R421").Value = 4.61
$ws.Range("S421").Value = 23
$ws.Range("T421").Value = 6
$ws.Range("U421").Value = 15
$ws.Range("V421").Value = 11

# Row 422: Karahali Souaré (right forward)
$ws.Range("A422").Value = 'Entrainement'
$ws.Range("B422").Value = 45903
$ws.Range("C422").Value = 'Global'
$ws.Range("D422").Value = 'J-3'
$ws.Range("E422").Value = 'Karahali Souaré'
$ws.Range("F422").Value = 'right forward'
$ws.Range("G422").Value = '01:25:02'
$ws.Range("H422").Value = 6.12
$ws.Range("I422").Value = 0.95
$ws.Range("J422").Value = 5.15
$ws.Range("K422").Value = 0.47
$ws.Range("L422").Value = 0.35
$ws.Range("M422").Value = 0.15
$ws.Range("N422").Value = 0
$ws.Range("O422").Value = 13
$ws.Range("P422").Value = 4.23
$ws.Range("Q422").Value = 29.22
$ws.Range("R422").Value = 5.2
$ws.Range("S422").Value = 28
$ws.Range("T422").Value = 8
$ws.Range("U422").Value = 20
$ws.Range("V422").Value = 9

# Row 423: Emmanuel Valey (left forward)
$ws.Range("A423").Value = 'Entrainement'
$ws.Range("B423").Value = 45903
$ws.Range("C423").Value = 'Global'
$ws.Range("D423").Value = 'J-3'
$ws.Range("E423").Value = 'Emmanuel Valey'
$ws.Range("F423").Value = 'left forward'
$ws.Range("G423").Value = '01:25:25'
$ws.Range("H423").Value = 6.95
$ws.Range("I423").Value = 1.18
$ws.Range("J423").Value = 5.75
$ws.Range("K423").Value = 0.65
$ws.Range("L423").Value = 0.38
$ws.Range("M423").Value = 0.16
$ws.Range("N423").Value = 0
$ws.Range("O423").Value = 10
$ws.Range("P423").Value = 4.8
$ws.Range("Q423").Value = 28.75
$ws.Range("R423").Value = 4.52
$ws.Range("S423").Value = 25
$ws.Range("T423").Value = 6
$ws.Range("U423").Value = 17
$ws.Range("V423").Value = 5

# Row 424: Naim Dhib (center midfield)
$ws.Range("A424").Value = 'Entrainement'
$ws.Range("B424").Value = 45903
$ws.Range("C424").Value = 'Global'
$ws.Range("D424").Value = 'J-3'
$ws.Range("E424").Value = 'Naim Dhib'
$ws.Range("F424").Value = 'center midfield'
$ws.Range("G424").Value = '01:25:09'
$ws.Range("H424").Value = 7.77
$ws.Range("I424").Value = 1.24
$ws.Range("J424").Value = 6.51
$ws.Range("K424").Value = 0.83
$ws.Range("L424").Value = 0.35
$ws.Range("M424").Value = 0.08
$ws.Range("N424").Value = 0
$ws.Range("O424").Value = 8
$ws.Range("P424").Value = 5.42
$ws.Range("Q424").Value = 29.22
$ws.Range("R424").Value = 4.51
$ws.Range("S424").Value = 29
$ws.Range("T424").Value = 6
$ws.Range("U424").Value = 17
$ws.Range("V424").Value = 5

# Row 425: Ilan Ihaddadene (center midfield)
$ws.Range("A425").Value = 'Entrainement'
$ws.Range("B425").Value = 45903
$ws.Range("C425").Value = 'Global'
$ws.Range("D425").Value = 'J-3'
$ws.Range("E425").Value = 'Ilan Ihaddadene'
$ws.Range("F425").Value = 'center midfield'
$ws.Range("G425").Value = '01:28:04'
$ws.Range("H425").Value = 7.43
$ws.Range("I425").Value = 1.15
$ws.Range("J425").Value = 6.27
$ws.Range("K425").Value = 0.76
$ws.Range("L425").Value = 0.31
$ws.Range("M425").Value = 0.09
$ws.Range("N425").Value = 0
$ws.Range("O425").Value = 9
$ws.Range("P425").Value = 4.99
$ws.Range("Q425").Value = 28.52
$ws.Range("R425").Value = 4.37
$ws.Range("S425").Value = 29
$ws.Range("T425").Value = 2
$ws.Range("U425").Value = 13
$ws.Range("V425").Value = 0

# Row 426: Yoan Zouma (center back)
$ws.Range("A426").Value = 'Entrainement'
$ws.Range("B426").Value = 45903
$ws.Range("C426").Value = 'Global'
$ws.Range("D426").Value = 'J-3'
$ws.Range("E426").Value = 'Yoan Zouma'
$ws.Range("F426").Value = 'center back'
$ws.Range("G426").Value = '01:25:16'
$ws.Range("H426").Value = 5.66
$ws.Range("I426").Value = 0.58
$ws.Range("J426").Value = 5.07
$ws.Range("K426").Value = 0.44
$ws.Range("L426").Value = 0.15
$ws.Range("M426").Value = 0
$ws.Range("N426").Value = 0
$ws.Range("O426").Value = 0
$ws.Range("P426").Value = 3.87
$ws.Range("Q426").Value = 24.8
$ws.Range("R426").Value = 4.96
$ws.Range("S426").Value = 19
$ws.Range("T426").Value = 2
$ws.Range("U426").Value = 15
$ws.Range("V426").Value = 4

# Row 427: Hedi Nasri (right back)
$ws.Range("A427").Value = 'Entrainement'
$ws.Range("B427").Value = 45903
$ws.Range("C427").Value = 'Global'
$ws.Range("D427").Value = 'J-3'
$ws.Range("E427").Value = 'Hedi Nasri'
$ws.Range("F427").Value = 'right back'
$ws.Range("G427").Value = '01:25:24'
$ws.Range("H427").Value = 6.45
$ws.Range("I427").Value = 0.85
$ws.Range("J427").Value = 5.59
$ws.Range("K427").Value = 0.44
$ws.Range("L427").Value = 0.29
$ws.Range("M427").Value = 0.13
$ws.Range("N427").Value = 0
$ws.Range("O427").Value = 11
$ws.Range("P427").Value = 4.48
$ws.Range("Q427").Value = 27.37
$ws.Range("R427").Value = 4.53
$ws.Range("S427").Value = 26
$ws.Range("T427").Value = 5
$ws.Range("U427").Value = 21
$ws.Range("V427").Value = 6

# Row 428: Sofiane Belle (left forward)
$ws.Range("A428").Value = 'Entrainement'
$ws.Range("B428").Value = 45903
$ws.Range("C428").Value = 'Global'
$ws.Range("D428").Value = 'J-3'
$ws.Range("E428").Value = 'Sofiane Belle'
$ws.Range("F428").Value = 'left forward'
$ws.Range("G428").Value = '01:26:55'
$ws.Range("H428").Value = 6.73
$ws.Range("I428").Value = 1.01
$ws.Range("J428").Value = 5.7
$ws.Range("K428").Value = 0.62
$ws.Range("L428").Value = 0.29
$ws.Range("M428").Value = 0.11
$ws.Range("N428").Value = 0
$ws.Range("O428").Value = 11
$ws.Range("P428").Value = 4.52
$ws.Range("Q428").Value = 27.24
$ws.Range("R428").Value = 4.11
$ws.Range("S428").Value = 27
$ws.Range("T428").Value = 1
$ws.Range("U428").Value = 17
$ws.Range("V428").Value = 3

# Row 429: Levy Ndoutoume (left back)
$ws.Range("A429").Value = 'Entrainement'
$ws.Range("B429").Value = 45903
$ws.Range("C429").Value = 'Global'
$ws.Range("D429").Value = 'J-3'
$ws.Range("E429").Value = 'Levy Ndoutoume'
$ws.Range("F429").Value = 'left back'
$ws.Range("G429").Value = '01:25:40'
$ws.Range("H429").Value = 6.98
$ws.Range("I429").Value = 0.86
$ws.Range("J429").Value = 6.11
$ws.Range("K429").Value = 0.6
$ws.Range("L429").Value = 0.26
$ws.Range("M429").Value = 0.01
$ws.Range("N429").Value = 0
$ws.Range("O429").Value = 1
$ws.Range("P429").Value = 4.82
$ws.Range("Q429").Value = 25.87
$ws.Range("R429").Value = 4.33
$ws.Range("S429").Value = 28
$ws.Range("T429").Value = 2
$ws.Range("U429").Value = 7
$ws.Range("V429").Value = 0

# Leave the view the way the author left it: scrolled down near the bottom
# of the sheet, with the next empty row in column D selected.
$excel.ActiveWindow.ScrollRow = 394
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D432").Select()

